$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Rectangle 5" (id=6) -> shape index 4
# target EMU: off x=4037902 y=3055896, ext cx=149887 cy=1394990
$shp1 = $s.Shapes.Item(4)
$shp1.Left = 317.945068359375
$shp1.Top = 240.62173461914062
$shp1.Width = 11.80212688446045
$shp1.Height = 109.84173583984375

# "Rectangle 17" (id=18) -> shape index 7
# target EMU: off x=5791122 y=3183991, ext cx=149259 cy=727122
$shp2 = $s.Shapes.Item(7)
$shp2.Left = 455.9938659667969
$shp2.Top = 250.7079620361328
$shp2.Width = 11.752677917480469
$shp2.Height = 57.25370407104492
